$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.305005550384521
$ws.Range("E2").Value = 851.766690651797
$ws.Range("F2").Value = 0.04461639847073016
$ws.Range("G2").Value = 0.03243901363710654
$ws.Range("H2").Value = 0.03044606920912101
$ws.Range("I2").Value = 0.02660201244329346
$ws.Range("J2").Value = 0.02637099712875958
$ws.Range("K2").Value = 0.02357672271373299
$ws.Range("L2").Value = 0.0231007271595054
$ws.Range("M2").Value = 0.02120178236391151
$ws.Range("N2").Value = 0.01989255730891966
$ws.Range("O2").Value = 0.0194553209747155
$ws.Range("P2").Value = 0.01780105741251484
$ws.Range("Q2").Value = 0.01780105741251484
$ws.Range("R2").Value = 0.01780105741251484
$ws.Range("S2").Value = 0.01750600542627236
$ws.Range("T2").Value = 0.0173670899626962
$ws.Range("U2").Value = 0.01710828239560704
$ws.Range("V2").Value = 0.01690096923101732
$ws.Range("W2").Value = 0.01680649263564861
$ws.Range("X2").Value = 0.01669356384599649
$ws.Range("Y2").Value = 0.01660363919399214

$ws.Range("C3").Value = 1.368995428085327
$ws.Range("E3").Value = 874.6324544991239
$ws.Range("F3").Value = 0.0434544730670968
$ws.Range("G3").Value = 0.03527205785988447
$ws.Range("H3").Value = 0.02878227267502709
$ws.Range("I3").Value = 0.02803011120699287
$ws.Range("J3").Value = 0.02409476295666809
$ws.Range("K3").Value = 0.02326886552795386
$ws.Range("L3").Value = 0.02151150724803208
$ws.Range("M3").Value = 0.02090304371883578
$ws.Range("N3").Value = 0.01999101683314505
$ws.Range("O3").Value = 0.01900762453532197
$ws.Range("P3").Value = 0.01857649619336955
$ws.Range("Q3").Value = 0.0182851854727198
$ws.Range("R3").Value = 0.01796576178477818
$ws.Range("S3").Value = 0.01787968668837649
$ws.Range("T3").Value = 0.01759955508466409
$ws.Range("U3").Value = 0.01728912624119711
$ws.Range("V3").Value = 0.01728912624119711
$ws.Range("W3").Value = 0.01718925739406081
$ws.Range("X3").Value = 0.01713236260226915
$ws.Range("Y3").Value = 0.01704936558477824

$ws.Range("C4").Value = 1.421000242233276
$ws.Range("E4").Value = 902.402881040025
$ws.Range("F4").Value = 0.04465459779664575
$ws.Range("G4").Value = 0.03380801703096048
$ws.Range("H4").Value = 0.02998254997601633
$ws.Range("I4").Value = 0.02727520660432462
$ws.Range("J4").Value = 0.02415178569633405
$ws.Range("K4").Value = 0.02215362233219
$ws.Range("L4").Value = 0.02104397999116079
$ws.Range("M4").Value = 0.01999239727170278
$ws.Range("N4").Value = 0.0198326490537921
$ws.Range("O4").Value = 0.01918579663232639
$ws.Range("P4").Value = 0.01918579663232639
$ws.Range("Q4").Value = 0.01907014613821567
$ws.Range("R4").Value = 0.01890915854437544
$ws.Range("S4").Value = 0.01861137524560301
$ws.Range("T4").Value = 0.01843069277144846
$ws.Range("U4").Value = 0.01790849064001159
$ws.Range("V4").Value = 0.01790849064001159
$ws.Range("W4").Value = 0.01769675393069265
$ws.Range("X4").Value = 0.01769675393069265
$ws.Range("Y4").Value = 0.01759069943547807

$ws.Range("C5").Value = 1.183997631072998
$ws.Range("E5").Value = 898.0682056365113
$ws.Range("F5").Value = 0.04430624886943869
$ws.Range("G5").Value = 0.03566807881646188
$ws.Range("H5").Value = 0.03042019152077201
$ws.Range("I5").Value = 0.02649481746006507
$ws.Range("J5").Value = 0.02391993207340458
$ws.Range("K5").Value = 0.02288020828371191
$ws.Range("L5").Value = 0.02193284670466573
$ws.Range("M5").Value = 0.02130802232331508
$ws.Range("N5").Value = 0.01977241440798943
$ws.Range("O5").Value = 0.01977241440798943
$ws.Range("P5").Value = 0.01945730229491135
$ws.Range("Q5").Value = 0.01889110099806826
$ws.Range("R5").Value = 0.01858466795853941
$ws.Range("S5").Value = 0.01854770680087788
$ws.Range("T5").Value = 0.01795076609802849
$ws.Range("U5").Value = 0.01777982878205156
$ws.Range("V5").Value = 0.01774378264622854
$ws.Range("W5").Value = 0.01755462798429056
$ws.Range("X5").Value = 0.01755462798429056
$ws.Range("Y5").Value = 0.01750620283891834

$ws.Range("C6").Value = 1.159997224807739
$ws.Range("E6").Value = 910.6015182520932
$ws.Range("F6").Value = 0.04387034731720904
$ws.Range("G6").Value = 0.0347910633317176
$ws.Range("H6").Value = 0.03034266669368169
$ws.Range("I6").Value = 0.02673611234633352
$ws.Range("J6").Value = 0.02333066885740094
$ws.Range("K6").Value = 0.02325884745652523
$ws.Range("L6").Value = 0.02206010123471652
$ws.Range("M6").Value = 0.02130822304928019
$ws.Range("N6").Value = 0.01988971820390992
$ws.Range("O6").Value = 0.01937983300356181
$ws.Range("P6").Value = 0.01925775489101959
$ws.Range("Q6").Value = 0.01887335930149878
$ws.Range("R6").Value = 0.01855828362348683
$ws.Range("S6").Value = 0.0183190509592909
$ws.Range("T6").Value = 0.0183104041948123
$ws.Range("U6").Value = 0.0181194663485884
$ws.Range("V6").Value = 0.01801038822055045
$ws.Range("W6").Value = 0.01783987082863542
$ws.Range("X6").Value = 0.01775051692499207
$ws.Range("Y6").Value = 0.01775051692499207

$ws.Range("C7").Value = 1.232996940612793
$ws.Range("E7").Value = 917.1301915141448
$ws.Range("F7").Value = 0.0444722041294939
$ws.Range("G7").Value = 0.03534558046979889
$ws.Range("H7").Value = 0.0319354870543627
$ws.Range("I7").Value = 0.02775507026129982
$ws.Range("J7").Value = 0.02595038433668952
$ws.Range("K7").Value = 0.02432482557883855
$ws.Range("L7").Value = 0.02237790264778471
$ws.Range("M7").Value = 0.02138517738692357
$ws.Range("N7").Value = 0.01976942256292904
$ws.Range("O7").Value = 0.01976942256292904
$ws.Range("P7").Value = 0.01946107285603292
$ws.Range("Q7").Value = 0.01899327317556576
$ws.Range("R7").Value = 0.01887279591297758
$ws.Range("S7").Value = 0.0187237979903564
$ws.Range("T7").Value = 0.01842200195293347
$ws.Range("U7").Value = 0.01841106564940834
$ws.Range("V7").Value = 0.0181912081633987
$ws.Range("W7").Value = 0.01808447626887423
$ws.Range("X7").Value = 0.01799062854291953
$ws.Range("Y7").Value = 0.01787778151099697

$ws.Range("C8").Value = 1.242002725601196
$ws.Range("E8").Value = 919.9091499951992
$ws.Range("F8").Value = 0.04451769468631686
$ws.Range("G8").Value = 0.03197527816349432
$ws.Range("H8").Value = 0.02995123243135885
$ws.Range("I8").Value = 0.02675130398150486
$ws.Range("J8").Value = 0.0259565787887662
$ws.Range("K8").Value = 0.02388319101200349
$ws.Range("L8").Value = 0.02325506421341361
$ws.Range("M8").Value = 0.02234585244204017
$ws.Range("N8").Value = 0.02139828456141693
$ws.Range("O8").Value = 0.01999419294145042
$ws.Range("P8").Value = 0.01977988618421092
$ws.Range("Q8").Value = 0.01908262560858954
$ws.Range("R8").Value = 0.01895793802472854
$ws.Range("S8").Value = 0.01895483573844277
$ws.Range("T8").Value = 0.0183596799989525
$ws.Range("U8").Value = 0.0183596799989525
$ws.Range("V8").Value = 0.0182768383029465
$ws.Range("W8").Value = 0.0181305727547759
$ws.Range("X8").Value = 0.01804585080435139
$ws.Range("Y8").Value = 0.01793195224162181

$ws.Range("C9").Value = 1.049999952316284
$ws.Range("E9").Value = 927.1509021990787
$ws.Range("F9").Value = 0.04509667442095978
$ws.Range("G9").Value = 0.03665096794071954
$ws.Range("H9").Value = 0.03212922397404209
$ws.Range("I9").Value = 0.02916262826187395
$ws.Range("J9").Value = 0.02662232066685439
$ws.Range("K9").Value = 0.02411898150173592
$ws.Range("L9").Value = 0.02190544906941409
$ws.Range("M9").Value = 0.02080764454452656
$ws.Range("N9").Value = 0.02076528758226095
$ws.Range("O9").Value = 0.02028608273660186
$ws.Range("P9").Value = 0.01923440888903961
$ws.Range("Q9").Value = 0.01923440888903961
$ws.Range("R9").Value = 0.01923440888903961
$ws.Range("S9").Value = 0.01908696654225295
$ws.Range("T9").Value = 0.0187454289824599
$ws.Range("U9").Value = 0.01866823189172384
$ws.Range("V9").Value = 0.01845686491549348
$ws.Range("W9").Value = 0.01822911987503026
$ws.Range("X9").Value = 0.01813832731706721
$ws.Range("Y9").Value = 0.01807311700193135

$ws.Range("C10").Value = 1.318007707595825
$ws.Range("E10").Value = 911.4635020143451
$ws.Range("F10").Value = 0.04357593814982216
$ws.Range("G10").Value = 0.03553456830912664
$ws.Range("H10").Value = 0.02878567404603457
$ws.Range("I10").Value = 0.02781751323179574
$ws.Range("J10").Value = 0.02693320290480431
$ws.Range("K10").Value = 0.02571314743971508
$ws.Range("L10").Value = 0.0243876529389581
$ws.Range("M10").Value = 0.0230125474293855
$ws.Range("N10").Value = 0.02265777527651826
$ws.Range("O10").Value = 0.02111694977926074
$ws.Range("P10").Value = 0.02060718745155351
$ws.Range("Q10").Value = 0.02032658530770004
$ws.Range("R10").Value = 0.01943961474866919
$ws.Range("S10").Value = 0.01922566047956718
$ws.Range("T10").Value = 0.01857751986911808
$ws.Range("U10").Value = 0.0182734113453769
$ws.Range("V10").Value = 0.01805440101080355
$ws.Range("W10").Value = 0.01791740280093362
$ws.Range("X10").Value = 0.01785388520988454
$ws.Range("Y10").Value = 0.01776731972737514

$ws.Range("C11").Value = 1.218000888824463
$ws.Range("E11").Value = 934.2241717391644
$ws.Range("F11").Value = 0.04438798055301021
$ws.Range("G11").Value = 0.03463524761627822
$ws.Range("H11").Value = 0.02948010945279311
$ws.Range("I11").Value = 0.02720146160545173
$ws.Range("J11").Value = 0.02632371822122969
$ws.Range("K11").Value = 0.02476211497892341
$ws.Range("L11").Value = 0.02279600589706349
$ws.Range("M11").Value = 0.02241868421440429
$ws.Range("N11").Value = 0.02059223813543031
$ws.Range("O11").Value = 0.02056772173536505
$ws.Range("P11").Value = 0.0200340398720046
$ws.Range("Q11").Value = 0.01889120088728964
$ws.Range("R11").Value = 0.01889120088728964
$ws.Range("S11").Value = 0.01889120088728964
$ws.Range("T11").Value = 0.01889120088728964
$ws.Range("U11").Value = 0.01872988689077452
$ws.Range("V11").Value = 0.01855733593714901
$ws.Range("W11").Value = 0.01835534964016067
$ws.Range("X11").Value = 0.01832567291430599
$ws.Range("Y11").Value = 0.01821099749978878
